$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.529999999999999
$ws.Range("A9").Value = -21.721
$ws.Range("D12").Value = -7.269
$ws.Range("E15").Value = 16.158
$ws.Range("A18").Value = -22.289
$ws.Range("A20").Value = -20.2
$ws.Range("D26").Value = -8.062000000000001
$ws.Range("A27").Value = -22.066
$ws.Range("D27").Value = -8.815999999999999
$ws.Range("D29").Value = -7.292
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.796000000000001
$ws.Range("E38").Value = 16.723
$ws.Range("E44").Value = 16.788
$ws.Range("D51").Value = -8.401
$ws.Range("E51").Value = 16.634
$ws.Range("D55").Value = -7.918000000000001
$ws.Range("E57").Value = 16.532
$ws.Range("E63").Value = 17.601
$ws.Range("A69").Value = -21.627
$ws.Range("D69").Value = -7.220999999999999
$ws.Range("D70").Value = -7.175999999999999
$ws.Range("E70").Value = 17.585
$ws.Range("A76").Value = -20.043
$ws.Range("A82").Value = -21.993
$ws.Range("D83").Value = -8.373000000000001
$ws.Range("E99").Value = 16.583
$ws.Range("D102").Value = -7.833
